$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume snapshot values (and the Chainlink/WrappedEther
# row swap) to match the refreshed data pulled by the scheduled GitHub
# Actions job. Price cells whose new text parses as a plain number are
# written with a leading apostrophe so Excel keeps them as literal text
# (matching the original inlineStr cells) instead of silently converting
# them to numeric values.

$ws.Range("D2").Value = "34.504.36"
$ws.Range("E2").Value = "  -2.87%  "
$ws.Range("D3").Value = "1.799.66"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'229.07"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "'0.613"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "'38.88"
$ws.Range("E8").Value = "  -10.42%  "
$ws.Range("E9").Value = "  +2.61%  "
$ws.Range("D10").Value = "'0.0676"
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("D11").Value = "'0.0988"
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("D12").Value = "2.059.26"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.792.43"
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.04"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("E16").Value = "  -3.89%  "
$ws.Range("D17").Value = "34.498.66"
$ws.Range("E17").Value = "  -2.86%  "
$ws.Range("E18").Value = "  -2.13%  "
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("D20").Value = "'238.45"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").Value = "'11.74"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").Value = "'4.67"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("D25").Value = "'172.58"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -4.17%  "
$ws.Range("D27").Value = "'17.12"
$ws.Range("E27").Value = "  -4.01%  "
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("E29").Value = "  -4.21%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "'4.00"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").Value = "'3.88"
$ws.Range("E33").Value = "  -5.40%  "
$ws.Range("E34").Value = "  +8.55%  "
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "'90.65"
$ws.Range("E37").Value = "  -4.81%  "
$ws.Range("E38").Value = "  +4.16%  "
$ws.Range("D39").Value = "1.311.81"
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("D41").Value = "'0.955"
$ws.Range("E41").Value = "  -5.74%  "
$ws.Range("D42").Value = "'2.42"
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("D43").Value = "'14.18"
$ws.Range("E43").Value = "  -8.22%  "
$ws.Range("E44").Value = "  -10.24%  "
$ws.Range("E45").Value = "  -4.19%  "
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").Value = "1.981.15"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("E50").Value = "  +3.59%  "
$ws.Range("D51").Value = "'97.34"
$ws.Range("E51").Value = "  -5.12%  "
